$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 24.35000000000037
$ws.Range("H2").Value = [double]"1.482768647245618e-16"
$ws.Range("K2").Value = 57.34701107123192
$ws.Range("L2").Value = "[53.49745296139423, 61.196569181069606]"
$ws.Range("O2").Value = 1.616395018964117
$ws.Range("P2").Value = "[1.5409213215805782, 1.6918687163476562]"
$ws.Range("S2").Value = 54.88419412935608
$ws.Range("T2").Value = "[52.31785161850519, 57.45053664020698]"
$ws.Range("W2").Value = 18.08578578578606
$ws.Range("X2").Value = 17.79329329329356
$ws.Range("Y2").Value = 18.37827827827856

$ws.Range("E3").Value = 25.24000000000051
$ws.Range("H3").Value = [double]"1.482768647245618e-16"
$ws.Range("K3").Value = 57.74272860300625
$ws.Range("L3").Value = "[52.71443812024369, 62.771019085768806]"
$ws.Range("O3").Value = -2.855421551010543
$ws.Range("P3").Value = "[-2.943474197958005, -2.7673689040630807]"
$ws.Range("S3").Value = 54.30462137494557
$ws.Range("T3").Value = "[51.42499756581498, 57.18424518407616]"
$ws.Range("W3").Value = 11.47043043043066
$ws.Range("X3").Value = 11.11671671671694
$ws.Range("Y3").Value = 11.82414414414438
